$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Archive an additional steady-state model run: add row 18 with the new run id
$ws.Range("A18").Value = 20211207

# Update the sheet's active selection/scroll position to the newly added cell
$ws.Range("B18").Select() | Out-Null
